# Daily attendance processing - reorder the "Recorded By" (column G) names
# for rows where the recorder list is out of the canonical order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact before -> after text replacements observed for the "Recorded By" column.
$map = @{
    "System, backup@backdoor.com, system" = "system, System, backup@backdoor.com";
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com";
    "admin@admin.com, System"             = "System, admin@admin.com";
    "dnasr281@gmail.com, admin@admin.com" = "admin@admin.com, dnasr281@gmail.com";
}

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 2; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($map.ContainsKey($val)) {
        $cell.Value2 = $map[$val]
    }
}
